$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 2.9
$ws.Range("L2").Value = 4
$ws.Range("W2").Value = 5.5
$ws.Range("AH2").Value = 12
$ws.Range("AJ2").Value = 34
$ws.Range("AQ2").Value = 67
$ws.Range("I4").Value = 3.2
$ws.Range("K4").Value = 1.95
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("AY4").Value = 34
$ws.Range("BB4").Value = 301
$ws.Range("G5").Value = 1.62
$ws.Range("I5").Value = 6.25
$ws.Range("N5").Value = 7.5
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("X5").Value = 6.5
$ws.Range("AC5").Value = 7.5
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 81
$ws.Range("AG5").Value = 13
$ws.Range("AO5").Value = 8.5
$ws.Range("AV5").Value = 81
$ws.Range("G7").Value = 2.15
$ws.Range("I7").Value = 3.1
$ws.Range("J7").Value = 2.75
$ws.Range("L7").Value = 3.5
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.38
$ws.Range("Z7").Value = 21
$ws.Range("AF7").Value = 34
$ws.Range("AL7").Value = 26
$ws.Range("AN7").Value = 4.5
$ws.Range("AU7").Value = 7
$ws.Range("AY7").Value = 21
$ws.Range("BA7").Value = 51
$ws.Range("G8").Value = 2.25
$ws.Range("I8").Value = 2.9
$ws.Range("J8").Value = 2.88
$ws.Range("L8").Value = 3.5
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.1
$ws.Range("X8").Value = 12
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 17
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 17
$ws.Range("AI8").Value = 11
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 19
$ws.Range("AZ8").Value = 51
$ws.Range("BA8").Value = 67
$ws.Range("G13").Value = 1.62
$ws.Range("I13").Value = 6.5
$ws.Range("L13").Value = 7.5
$ws.Range("W13").Value = 4.33
$ws.Range("AC13").Value = 5.5
$ws.Range("AI13").Value = 23
$ws.Range("AK13").Value = 67
$ws.Range("AU13").Value = 12
$ws.Range("M16").Value = 1.11
$ws.Range("N16").Value = 6.5
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 3.25
$ws.Range("I17").Value = 2.38
$ws.Range("J17").Value = 3.75
$ws.Range("L17").Value = 3.1
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.8
$ws.Range("W17").Value = 8.5
$ws.Range("X17").Value = 15
$ws.Range("Z17").Value = 34
$ws.Range("AB17").Value = 41
$ws.Range("AC17").Value = 8.5
$ws.Range("AD17").Value = 6.5
$ws.Range("AE17").Value = 17
$ws.Range("AG17").Value = 7
$ws.Range("AH17").Value = 11
$ws.Range("AN17").Value = 5
$ws.Range("AS17").Value = 251
$ws.Range("AV17").Value = 67
$ws.Range("AW17").Value = 4.33
$ws.Range("AX17").Value = 13
$ws.Range("BA17").Value = 67
$ws.Range("M20").Value = 1.06
